$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Strip the "Heading2" paragraph style from the four section headings
#    (Introduction, Evidence and Analysis, Counterarguments, Conclusion).
#    Re-applying "Normal" is the COM-level equivalent of clearing the style
#    in the Word UI (Styles gallery -> Normal).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 2") {
        $p.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the inline citation markers for the new "Ref-xxxxxx" placeholders.
#    Every replacement is scoped to the specific paragraph's Range so that
#    identical citation text in other paragraphs (e.g. "(Sun and Sun)" or
#    "(Addo and Fang)" appear more than once across the document, each time
#    mapping to a different Ref code) is left untouched.
# ---------------------------------------------------------------------------

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $rng = $d.Paragraphs($paraIndex).Range
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $newText, 2)
}

# Paragraph 8 - "Recent studies have showcased ..."
Replace-InParagraph 8 "(Franceschini et al.)" "(Ref-u381640)"

# Paragraph 10 - "However, the association between video games ..."
Replace-InParagraph 10 "(Verheijen et al.)" "(Ref-f669144)"
Replace-InParagraph 10 "(Sun and Sun)" "(Ref-f669144)"

# Paragraph 12 - "Additionally, video games offer a complex interplay ..."
Replace-InParagraph 12 "(López-Fernández and Mezquita)" "(Ref-f203639)"

# Paragraph 16 - "Despite the concerns raised about video games ..."
Replace-InParagraph 16 "(Sun and Sun)" "(Ref-s336207)"
Replace-InParagraph 16 "(Addo and Fang)" "(Ref-s336207)"

# Paragraph 18 - "Moreover, expert opinions further bolster ..."
Replace-InParagraph 18 "Addo and Fang," "Ref-A1B2C3,"
Replace-InParagraph 18 "(Addo and Fang)" "(Ref-s486805)"
Replace-InParagraph 18 "Verheijen et al.," "Ref-G7H3JS,"
Replace-InParagraph 18 "(Verheijen et al.)" "(Ref-s486805)"
